# SwaadSutra_Daily_2026-01-13.xlsx update
# New order #13 (Ketki, Wheat Chapati x1) came in at 2026-01-13 22:51.
# It is inserted as the new top row of the "Daily Orders" log (the sheet
# is sorted newest-first), pushing every existing order down by one row.
# The Summary and Items Breakdown aggregate sheets are updated to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Daily Orders — insert the new order at the top (row 2)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Daily Orders")

$ws.Rows.Item(2).Insert()

$ws.Cells.Item(2, 1).Value = 13
$ws.Cells.Item(2, 2).Value = "2026-01-13 22:51"
$ws.Cells.Item(2, 3).Value = "Ketki"

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "1608"
$ws.Cells.Item(2, 4).Style = "Normal"

$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "3159135521"
$ws.Cells.Item(2, 5).Style = "Normal"

$ws.Cells.Item(2, 6).Value = "Wheat Chapati x1"
$ws.Cells.Item(2, 7).Value = 15
$ws.Cells.Item(2, 8).Value = "NEW"
$ws.Cells.Item(2, 9).Value = "PENDING"

$ws.Cells.Item(2, 10).NumberFormat = "@"
$ws.Cells.Item(2, 10).Value = "2026-01-14"
$ws.Cells.Item(2, 10).Style = "Normal"

$ws.Cells.Item(2, 11).NumberFormat = "@"
$ws.Cells.Item(2, 11).Value = "16:51"
$ws.Cells.Item(2, 11).Style = "Normal"

# Notes / Cancel Reason / Feedback are blank for this order.
$ws.Cells.Item(2, 12).Value = ""
$ws.Cells.Item(2, 13).Value = ""
$ws.Cells.Item(2, 14).Value = ""

# ---------------------------------------------------------------
# 2) Summary — bump Total Orders / New / Total Revenue
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Cells.Item(2, 1).Value = 13   # Total Orders
$summary.Cells.Item(2, 2).Value = 11   # New
$summary.Cells.Item(2, 7).Value = 340  # Total Revenue

# ---------------------------------------------------------------
# 3) Items Breakdown — Wheat Chapati qty/revenue go up by 1 / 15
# ---------------------------------------------------------------
$items = $wb.Worksheets.Item("Items Breakdown")
$items.Cells.Item(2, 2).Value = 6   # Quantity Ordered
$items.Cells.Item(2, 3).Value = 90  # Revenue
